$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.295.45'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '1.908.76'
$ws.Range('E3').Value = '  +2.23%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.64'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5256'
$ws.Range('E7').Value = '  +3.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3789'
$ws.Range('E8').Value = '  +3.65%  '
$ws.Range('E9').Value = '  +1.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.31'
$ws.Range('E10').Value = '  +3.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9015'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08161'
$ws.Range('E12').Value = '  +9.05%  '
$ws.Range('D13').Value = '1.914.12'
$ws.Range('E13').Value = '  +2.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '95.39'
$ws.Range('E14').Value = '  +1.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.299'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008615'
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.51'
$ws.Range('E18').Value = '  +2.73%  '
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').Value = '27.344.81'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.074'
$ws.Range('E21').Value = '  +1.47%  '
$ws.Range('D22').Value = '2.153.09'
$ws.Range('E22').Value = '  +2.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.66'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.466'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.310'
$ws.Range('E25').Value = '  +10.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '146.19'
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('E27').Value = '  -1.79%  '
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '115.06'
$ws.Range('E29').Value = '  +1.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.000'
$ws.Range('E30').Value = '  +6.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.813'
$ws.Range('E31').Value = '  +3.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09238'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8067'
$ws.Range('E33').Value = '  +7.99%  '
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('E35').Value = '  +8.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.960'
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.359'
$ws.Range('E37').Value = '  +4.38%  '
$ws.Range('E38').Value = '  +2.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5741'
$ws.Range('E39').Value = '  +3.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01985'
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.078'
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.999'
$ws.Range('E42').Value = '  +4.59%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.64'
$ws.Range('E43').Value = '  +3.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.629'
$ws.Range('E44').Value = '  +0.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1518'
$ws.Range('E45').Value = '  +2.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4851'
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.21'
$ws.Range('E47').Value = '  +1.67%  '
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('E49').Value = '  +4.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.68'
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('E51').Value = '  +1.42%  '
